$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Fix the typo "Bando de dados:" -> "Banco de dados:" the way Word itself
# would: replace the single character "d" (4th letter) with "c", which is
# what causes Word to split the run and drop a fresh _GoBack bookmark at the
# point of the edit (and remove it from its previous location).
# ---------------------------------------------------------------------------

$full = $d.Content.Text
$idx = $full.IndexOf("Bando de dados:")

if ($idx -ge 0) {
    $typoPos = $idx + 3

    # Replace the "d" in "Bando" with "c" -> "Banco"
    $delRange = $d.Range($typoPos, $typoPos + 1)
    $delRange.Delete()

    $insRange = $d.Range($typoPos, $typoPos)
    $insRange.InsertAfter("c")

    # Force a run boundary right after "Ban" (before the newly typed "c")
    # using a scratch bookmark, so "Ban" and "c" end up as distinct runs --
    # exactly like Word leaves behind after a keystroke-level edit.
    $splitRange = $d.Range($typoPos, $typoPos)
    $d.Bookmarks.Add("__scratch_split", $splitRange) | Out-Null

    # Drop the real _GoBack bookmark right after the "c" (this also moves it
    # off of its old location elsewhere in the document, since a document can
    # only have a single bookmark with a given name).
    $goBackRange = $d.Range($typoPos + 1, $typoPos + 1)
    $d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

    # Remove the scratch bookmark now that both run boundaries are in place.
    $d.Bookmarks.Item("__scratch_split").Delete()
}
